# Refresh cryptos list price/volume data (and two coin-row reorderings)
# to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.055.57"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.005.73"
$ws.Range("E3").Value = "  -2.25%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.42%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.85"
$ws.Range("E5").Value = "  -2.22%  "

# Row 6: XRP
$ws.Range("E6").Value = "  -2.41%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.08%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.62"
$ws.Range("E8").Value = "  -4.29%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  -1.55%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").Value = "  +1.57%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -2.93%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.304.24"
$ws.Range("E12").Value = "  -2.10%  "

# Row 13: Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.10"
$ws.Range("E13").Value = "  -3.63%  "

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.25"
$ws.Range("E14").Value = "  -2.00%  "

# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.739"
$ws.Range("E15").Value = "  -2.41%  "

# Row 16: Polkadot
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.11"
$ws.Range("E16").Value = "  -3.21%  "

# Row 17: WrappedEther
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.008.88"
$ws.Range("E17").Value = "  -2.42%  "

# Row 18: WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.938.55"
$ws.Range("E18").Value = "  -1.01%  "

# Row 19: Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -0.10%  "

# Row 20: Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.71"
$ws.Range("E20").Value = "  -1.34%  "

# Row 21: ShibaInu
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0815"
$ws.Range("E21").Value = "  -0.98%  "

# Row 22: BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.17"
$ws.Range("E22").Value = "  -1.42%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.05%  "

# Row 24: Toncoin
$ws.Range("E24").Value = "  +1.17%  "

# Row 25: PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -6.40%  "

# Row 26: Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.00"
$ws.Range("E26").Value = "  -2.38%  "

# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  -7.56%  "

# Row 28: Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.124"
$ws.Range("E28").Value = "  -4.03%  "

# Row 29: EthereumClassic
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.64"
$ws.Range("E29").Value = "  -2.93%  "

# Row 30: ImmutableX
$ws.Range("E30").Value = "  -0.21%  "

# Row 31: Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("E31").Value = "  -3.59%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.50"
$ws.Range("E32").Value = "  -0.65%  "

# Row 33: Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0608"
$ws.Range("E33").Value = "  -2.17%  "

# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.40"
$ws.Range("E34").Value = "  -3.59%  "

# Row 35: LidoDAOToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.33"
$ws.Range("E35").Value = "  -5.92%  "

# Row 36: WEMIXToken
$ws.Range("E36").Value = "  +1.66%  "

# Row 37: BinanceUSD
$ws.Range("E37").Value = "  +0.42%  "

# Row 38: RenderToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.14"
$ws.Range("E38").Value = "  -4.22%  "

# Row 39: THORChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.32"
$ws.Range("E39").Value = "  +0.08%  "

# Row 40: Maker
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.478.60"
$ws.Range("E40").Value = "  +0.07%  "

# Row 41: VeChain
$ws.Range("E41").Value = "  -3.71%  "

# Row 42: Aave
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "94.54"
$ws.Range("E42").Value = "  -3.80%  "

# Row 43: InjectiveProtocol
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.47"
$ws.Range("E43").Value = "  -0.94%  "

# Row 44: Cronos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0918"
$ws.Range("E44").Value = "  -3.68%  "

# Row 45: HuobiToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").Value = "  -4.87%  "

# Row 46: TrustWalletToken
$ws.Range("E46").Value = "  -4.87%  "

# Row 47: FraxShare
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.17"
$ws.Range("E47").Value = "  -1.20%  "

# Row 48: ARBITRUM
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -2.35%  "

# Row 49: MXToken
$ws.Range("E49").Value = "  -1.09%  "

# Row 50: RocketPoolETH
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.194.10"
$ws.Range("E50").Value = "  -2.04%  "

# Row 51: MultiversX
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.36"
$ws.Range("E51").Value = "  -3.14%  "
